# Add new "hgt" (Hong Kong / shareholding) table rows sourced from tinysoft
# to the FactorID sheet: three new rows (23-25) describing table 132
# (截止日/date, 股数/shares_hold, 占总股本比例(%)/ratio_hold).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: 截止日 / 132001 / date / 132
$ws.Range("A23").Value = "截止日"
$ws.Range("B23").Value = 132001
$ws.Range("C23").Value = "date"
$ws.Range("D23").Value = 132

# New row 24: 股数 / 132002 / shares_hold / 132
$ws.Range("A24").Value = "股数"
$ws.Range("B24").Value = 132002
$ws.Range("C24").Value = "shares_hold"
$ws.Range("D24").Value = 132

# New row 25: 占总股本比例(%) / 132003 / ratio_hold / 132
$ws.Range("A25").Value = "占总股本比例(%)"
$ws.Range("B25").Value = 132003
$ws.Range("C25").Value = "ratio_hold"
$ws.Range("D25").Value = 132

# Column C for the English-name field uses the wrap-text/size-12 style seen
# on the other English-name cells in the table (e.g. C2) - copy that
# formatting onto the three new cells.
$ws.Range("C2").Copy()
$ws.Range("C23:C25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the author's final selection after entering the new data.
$ws.Range("D26").Select()
